# tablaTest.xlsx touch-up:
#  - the two previously-entered "Periodo de Imputacion" dates (E3:E5 had been
#    filled in, now they are blanked out again)
#  - the remembered selection moves from B12 back onto the table (B4)
#  - the sheet gets an explicit page setup (A4 / portrait) for printing
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Blank out E3:E5 (previously held date serials 44785/44784/44783); keep
# their existing style (the date number format) intact, just drop the values.
$ws.Range("E3").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("E5").ClearContents()

# Move the active selection/cursor to B4.
$ws.Range("B4").Select()

# Page setup: A4 paper, portrait orientation.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
